# Auto-generated Excel for date: 2025-01-14
# Shifts the 28 consecutive daily-schedule rows (column B) so that the
# first row starts on 2025-01-14 instead of 2024-12-26, keeping every
# date consecutive. All other columns/content are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "2025-01-14","2025-01-15","2025-01-16","2025-01-17","2025-01-18","2025-01-19","2025-01-20",
    "2025-01-21","2025-01-22","2025-01-23","2025-01-24","2025-01-25","2025-01-26","2025-01-27",
    "2025-01-28","2025-01-29","2025-01-30","2025-01-31","2025-02-01","2025-02-02","2025-02-03",
    "2025-02-04","2025-02-05","2025-02-06","2025-02-07","2025-02-08","2025-02-09","2025-02-10"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Range("B$row")
    # Force text storage so Excel doesn't auto-convert the "yyyy-mm-dd"
    # looking string into a real date serial number / date format.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    # Restore the default (unstyled) look so the cell keeps matching the
    # rest of the sheet, exactly like the original plain-text date cells.
    $cell.Style = "Normal"
}
